$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.508.89'
$ws.Cells.Item(2, 5).Value = '  -0.93%  '

$ws.Cells.Item(3, 4).Value = '2.294.30'
$ws.Cells.Item(3, 5).Value = '  -0.76%  '

$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  -0.04%  '

$ws.Cells.Item(5, 4).Value = '''300.69'
$ws.Cells.Item(5, 5).Value = '  -0.77%  '

$ws.Cells.Item(6, 4).Value = '''96.95'
$ws.Cells.Item(6, 5).Value = '  -3.23%  '

$ws.Cells.Item(7, 4).Value = '''0.499'
$ws.Cells.Item(7, 5).Value = '  -1.17%  '

$ws.Cells.Item(8, 5).Value = '  +0.04%  '

$ws.Cells.Item(9, 4).Value = '''0.494'
$ws.Cells.Item(9, 5).Value = '  -2.09%  '

$ws.Cells.Item(10, 4).Value = '''33.57'
$ws.Cells.Item(10, 5).Value = '  -3.78%  '

$ws.Cells.Item(11, 4).Value = '''0.0792'
$ws.Cells.Item(11, 5).Value = '  -0.02%  '

$ws.Cells.Item(12, 4).Value = '''48.45'
$ws.Cells.Item(12, 5).Value = '  -6.13%  '

$ws.Cells.Item(13, 5).Value = '  +2.16%  '

$ws.Cells.Item(14, 4).Value = '''16.13'
$ws.Cells.Item(14, 5).Value = '  +2.96%  '

$ws.Cells.Item(15, 4).Value = '''6.73'
$ws.Cells.Item(15, 5).Value = '  +0.01%  '

$ws.Cells.Item(16, 4).Value = '2.648.26'
$ws.Cells.Item(16, 5).Value = '  -0.80%  '

$ws.Cells.Item(17, 4).Value = '2.274.54'
$ws.Cells.Item(17, 5).Value = '  -2.27%  '

$ws.Cells.Item(18, 4).Value = '''0.794'
$ws.Cells.Item(18, 5).Value = '  -0.52%  '

$ws.Cells.Item(19, 4).Value = '42.424.50'
$ws.Cells.Item(19, 5).Value = '  -0.89%  '

$ws.Cells.Item(20, 4).Value = '''11.74'
$ws.Cells.Item(20, 5).Value = '  +0.47%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0897'
$ws.Cells.Item(21, 5).Value = '  -1.03%  '

$ws.Cells.Item(22, 4).Value = '''6.02'
$ws.Cells.Item(22, 5).Value = '  -0.67%  '

$ws.Cells.Item(23, 4).Value = '''66.68'
$ws.Cells.Item(23, 5).Value = '  -1.86%  '

$ws.Cells.Item(24, 4).Value = '''236.21'
$ws.Cells.Item(24, 5).Value = '  +0.04%  '

$ws.Cells.Item(25, 4).Value = '''1.97'
$ws.Cells.Item(25, 5).Value = '  +0.46%  '

$ws.Cells.Item(26, 2).Value = 'PancakeSwap'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(26, 4).Value = '''2.47'
$ws.Cells.Item(26, 5).Value = '  -1.97%  '

$ws.Cells.Item(27, 2).Value = 'Dai'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(27, 4).Value = '''1.00'
$ws.Cells.Item(27, 5).Value = '  -0.07%  '

$ws.Cells.Item(28, 4).Value = '''23.98'
$ws.Cells.Item(28, 5).Value = '  -3.56%  '

$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).Value = '''2.17'
$ws.Cells.Item(29, 5).Value = '  -0.03%  '

$ws.Cells.Item(30, 2).Value = 'Monero'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(30, 4).Value = '''167.52'
$ws.Cells.Item(30, 5).Value = '  +1.09%  '

$ws.Cells.Item(31, 4).Value = '''34.06'
$ws.Cells.Item(31, 5).Value = '  -1.47%  '

$ws.Cells.Item(32, 4).Value = '''9.18'
$ws.Cells.Item(32, 5).Value = '  +0.57%  '

$ws.Cells.Item(33, 4).Value = '''0.998'
$ws.Cells.Item(33, 5).Value = '  -0.11%  '

$ws.Cells.Item(34, 4).Value = '''4.72'
$ws.Cells.Item(34, 5).Value = '  +5.68%  '

$ws.Cells.Item(35, 4).Value = '''4.96'
$ws.Cells.Item(35, 5).Value = '  -1.25%  '

$ws.Cells.Item(36, 4).Value = '''16.98'
$ws.Cells.Item(36, 5).Value = '  +0.96%  '

$ws.Cells.Item(37, 4).Value = '''0.0697'
$ws.Cells.Item(37, 5).Value = '  -0.28%  '

$ws.Cells.Item(38, 5).Value = '  -3.20%  '

$ws.Cells.Item(39, 4).Value = '''2.81'
$ws.Cells.Item(39, 5).Value = '  -2.72%  '

$ws.Cells.Item(40, 4).Value = '''0.0994'
$ws.Cells.Item(40, 5).Value = '  -1.29%  '

$ws.Cells.Item(41, 4).Value = '''1.75'
$ws.Cells.Item(41, 5).Value = '  -3.13%  '

$ws.Cells.Item(42, 5).Value = '  -1.57%  '

$ws.Cells.Item(43, 5).Value = '  -9.52%  '

$ws.Cells.Item(44, 4).Value = '1.965.29'
$ws.Cells.Item(44, 5).Value = '  -0.36%  '

$ws.Cells.Item(45, 4).Value = '''0.0280'
$ws.Cells.Item(45, 5).Value = '  -0.30%  '

$ws.Cells.Item(46, 4).Value = '''17.72'
$ws.Cells.Item(46, 5).Value = '  -4.21%  '

$ws.Cells.Item(47, 4).Value = '''9.67'
$ws.Cells.Item(47, 5).Value = '  -5.74%  '

$ws.Cells.Item(48, 4).Value = '''2.82'
$ws.Cells.Item(48, 5).Value = '  -2.18%  '

$ws.Cells.Item(49, 4).Value = '2.512.36'
$ws.Cells.Item(49, 5).Value = '  -0.81%  '

$ws.Cells.Item(50, 4).Value = '''52.66'
$ws.Cells.Item(50, 5).Value = '  -5.42%  '

$ws.Cells.Item(51, 5).Value = '  -3.06%  '
